$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 401726.62
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 401726.62
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1205179.86
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1205515.86
$ws.Range("H39").Value = 2298.818
$ws.Range("I39").Value = 3411.5
$ws.Range("J39").Value = 963.6
$ws.Range("K39").Value = 10234.5
$ws.Range("L39").Value = 2890.8
$ws.Range("M39").Value = -9938.5
$ws.Range("N39").Value = -3482.8
$ws.Range("H62").Value = 5417.75
$ws.Range("I62").Value = 4383.5
$ws.Range("K62").Value = 4383.5
$ws.Range("M62").Value = -3759.5
$ws.Range("H65").Value = 5417.75
$ws.Range("I65").Value = 4383.5
$ws.Range("K65").Value = 21917.5
$ws.Range("M65").Value = -18797.5
$ws.Range("H80").Value = 1809.3334
$ws.Range("I80").Value = 1082.75
$ws.Range("J80").Value = 2390.6
$ws.Range("K80").Value = 3248.25
$ws.Range("L80").Value = 7171.799999999999
$ws.Range("M80").Value = -2250.25
$ws.Range("N80").Value = -9167.799999999999
$ws.Range("H83").Value = 1809.3334
$ws.Range("I83").Value = 1082.75
$ws.Range("J83").Value = 2390.6
$ws.Range("K83").Value = 9744.75
$ws.Range("L83").Value = 21515.4
$ws.Range("M83").Value = -4752.75
$ws.Range("N83").Value = -31499.4
$ws.Range("H86").Value = 2032.7693
$ws.Range("I86").Value = 1804
$ws.Range("K86").Value = 1804
$ws.Range("M86").Value = -681
$ws.Range("H89").Value = 2032.7693
$ws.Range("I89").Value = 1804
$ws.Range("K89").Value = 9020
$ws.Range("M89").Value = -3404
$ws.Range("H107").Value = 1011.8571
$ws.Range("I107").Value = 597.64703
$ws.Range("J107").Value = 2772.25
$ws.Range("K107").Value = 597.64703
$ws.Range("L107").Value = 2772.25
$ws.Range("M107").Value = 1322.35297
$ws.Range("N107").Value = -6612.25
$ws.Range("H133").Value = 80733
$ws.Range("J133").Value = 80733
$ws.Range("L133").Value = 80733
$ws.Range("N133").Value = -90853
$ws.Range("H138").Value = 2614.8108
$ws.Range("I138").Value = 1779.8667
$ws.Range("J138").Value = 3184.0908
$ws.Range("K138").Value = 5339.6001
$ws.Range("L138").Value = 9552.2724
$ws.Range("M138").Value = -199.6000999999997
$ws.Range("N138").Value = -19832.2724
$ws.Range("H141").Value = 2216.0908
$ws.Range("I141").Value = 2216.0908
$ws.Range("K141").Value = 6648.2724
$ws.Range("M141").Value = -1468.2724

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4018.1216
$ws.Range("I32").Value = 3702.1714
$ws.Range("J32").Value = 9547.25
$ws.Range("K32").Value = 3702.1714
$ws.Range("L32").Value = 9547.25
$ws.Range("M32").Value = -3415.1714
$ws.Range("N32").Value = -10121.25
$ws.Range("H43").Value = 10377
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10377
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10377
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -11003
$ws.Range("H61").Value = 1578.4615
$ws.Range("I61").Value = 1471.0217
$ws.Range("J61").Value = 2402.1667
$ws.Range("K61").Value = 1471.0217
$ws.Range("L61").Value = 2402.1667
$ws.Range("M61").Value = -1259.0217
$ws.Range("N61").Value = -2826.1667
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
$ws.Range("H122").Value = 2919.2144
$ws.Range("I122").Value = 2408.7222
$ws.Range("K122").Value = 7226.1666
$ws.Range("M122").Value = -4776.1666
$ws.Range("H132").Value = 7754.8486
$ws.Range("I132").Value = 8063.2
$ws.Range("J132").Value = 4671.3335
$ws.Range("K132").Value = 24189.6
$ws.Range("L132").Value = 14014.0005
$ws.Range("M132").Value = -21659.6
$ws.Range("N132").Value = -19074.0005
$ws.Range("H136").Value = 1578.4615
$ws.Range("I136").Value = 1471.0217
$ws.Range("J136").Value = 2402.1667
$ws.Range("K136").Value = 4413.0651
$ws.Range("L136").Value = 7206.500100000001
$ws.Range("M136").Value = -1863.0651
$ws.Range("N136").Value = -12306.5001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 154921.33
$ws.Range("J42").Value = 154921.33
$ws.Range("L42").Value = 154921.33
$ws.Range("N42").Value = -155577.33
$ws.Range("H105").Value = 2099.0908
$ws.Range("I105").Value = 2089.1
$ws.Range("K105").Value = 2089.1
$ws.Range("M105").Value = -342.0999999999999
$ws.Range("H134").Value = 177647.64
$ws.Range("I134").Value = 188713.19
$ws.Range("J134").Value = 599
$ws.Range("K134").Value = 566139.5700000001
$ws.Range("L134").Value = 1797
$ws.Range("M134").Value = -563604.5700000001
$ws.Range("N134").Value = -6867
$ws.Range("H141").Value = 87778
$ws.Range("J141").Value = 87778
$ws.Range("L141").Value = 87778
$ws.Range("N141").Value = -98138

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1453.7222
$ws.Range("I58").Value = 1489.3636
$ws.Range("J58").Value = 1397.7142
$ws.Range("K58").Value = 1489.3636
$ws.Range("L58").Value = 1397.7142
$ws.Range("M58").Value = -1286.3636
$ws.Range("N58").Value = -1803.7142
$ws.Range("H68").Value = 25000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 25000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H81").Value = 26298
$ws.Range("I81").Value = 26298
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 26298
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -25300
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 26298
$ws.Range("I84").Value = 26298
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 78894
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -73902
$ws.Range("N84").ClearContents()
$ws.Range("H136").Value = 1453.7222
$ws.Range("I136").Value = 1489.3636
$ws.Range("J136").Value = 1397.7142
$ws.Range("K136").Value = 4468.0908
$ws.Range("L136").Value = 4193.142599999999
$ws.Range("M136").Value = -1918.0908
$ws.Range("N136").Value = -9293.142599999999
$ws.Range("H141").Value = 330246.7
$ws.Range("J141").Value = 356940.78
$ws.Range("L141").Value = 356940.78
$ws.Range("N141").Value = -367300.78

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 534.0476
$ws.Range("I33").Value = 238.46153
$ws.Range("J33").Value = 1014.375
$ws.Range("K33").Value = 1430.76918
$ws.Range("L33").Value = 6086.25
$ws.Range("M33").Value = -1147.76918
$ws.Range("N33").Value = -6652.25

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3573.375
$ws.Range("I80").Value = 2599.5
$ws.Range("J80").Value = 6495
$ws.Range("K80").Value = 2599.5
$ws.Range("L80").Value = 6495
$ws.Range("M80").Value = -1601.5
$ws.Range("N80").Value = -8491
$ws.Range("H83").Value = 3573.375
$ws.Range("I83").Value = 2599.5
$ws.Range("J83").Value = 6495
$ws.Range("K83").Value = 12997.5
$ws.Range("L83").Value = 32475
$ws.Range("M83").Value = -8005.5
$ws.Range("N83").Value = -42459
$ws.Range("H95").Value = 47500
$ws.Range("I95").Value = 45000
$ws.Range("J95").Value = 50000
$ws.Range("K95").Value = 45000
$ws.Range("L95").Value = 50000
$ws.Range("M95").Value = -42254
$ws.Range("N95").Value = -55492
$ws.Range("H132").Value = 3685.617
$ws.Range("I132").Value = 2722.4102
$ws.Range("K132").Value = 8167.230599999999
$ws.Range("M132").Value = -5637.230599999999
$ws.Range("H136").Value = 49554.31
$ws.Range("J136").Value = 49554.31
$ws.Range("L136").Value = 148662.93
$ws.Range("N136").Value = -153762.93

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 12950.429
$ws.Range("I93").Value = 3145.577
$ws.Range("K93").Value = 3145.577
$ws.Range("M93").Value = -1897.577
$ws.Range("H109").Value = 55997.668
$ws.Range("J109").Value = 55997.668
$ws.Range("L109").Value = 55997.668
$ws.Range("N109").Value = -58771.668
$ws.Range("H140").Value = 62959.332
$ws.Range("J140").Value = 62959.332
$ws.Range("L140").Value = 62959.332
$ws.Range("N140").Value = -73319.33199999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30301
$ws.Range("J80").Value = 30301
$ws.Range("L80").Value = 30301
$ws.Range("N80").Value = -32297
$ws.Range("H83").Value = 30301
$ws.Range("J83").Value = 30301
$ws.Range("L83").Value = 90903
$ws.Range("N83").Value = -100887
$ws.Range("H95").Value = 29114.666
$ws.Range("J95").Value = 29114.666
$ws.Range("L95").Value = 29114.666
$ws.Range("N95").Value = -34606.666
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H104").Value = 29999.5
$ws.Range("J104").Value = 29999.5
$ws.Range("L104").Value = 29999.5
$ws.Range("N104").Value = -36987.5
$ws.Range("H105").Value = 40614.5
$ws.Range("J105").Value = 40614.5
$ws.Range("L105").Value = 40614.5
$ws.Range("N105").Value = -47602.5
$ws.Range("H113").Value = 423.7143
$ws.Range("I113").Value = 375.66666
$ws.Range("K113").Value = 1126.99998
$ws.Range("M113").Value = 1043.00002
$ws.Range("H136").Value = 1059.5186
$ws.Range("I136").Value = 1054.1154
$ws.Range("K136").Value = 3162.3462
$ws.Range("M136").Value = -612.3462
